# adding averages and more checks
#
# - PERIOD TO EXPIRE (col H, rows 3-13) drops by 8 days (re-run against a
#   later "as of" date).
# - LAST UPDATE (col I, rows 3-13) moves from 08-Sep-2025 to 16-Sep-2025.
# - Header styling: the title (A1) and column-header row (A2:K2) fonts are
#   unified into a single bold/white font (the old bold+14pt title font and
#   the plain bold header font collapse into one style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Data rows: PERIOD TO EXPIRE (H) and LAST UPDATE (I)
# ---------------------------------------------------------------------
$rows = 3..13

foreach ($r in $rows) {
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value() - 8
}

# LAST UPDATE (I3:I13) needs to become the literal text "16-Sep-2025" and
# stay text (these columns already hold "dd-mmm-yyyy" strings, not real
# dates). Assigning that string straight to .Value lets the host
# auto-recognise it as a date and store a date serial instead, so build it
# once as text on a scratch cell (forced text format) and paste-special
# the *values* into each target cell, which carries over the literal text
# without dragging the scratch cell's number format along.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "16-Sep-2025"
$scratch.Copy()

foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).PasteSpecial(-4163)   # xlPasteValues
}

$scratch.Clear()

# ---------------------------------------------------------------------
# 2) Header styling: merge the two bold fonts into one bold/white font
#    used by both the dashboard title and the column header row.
# ---------------------------------------------------------------------
$titleRange = $ws.Range("A1")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 11
$titleRange.Font.Color = 16777215

$headerRange = $ws.Range("A2:K2")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
